# Fixed up Myxicola and a few more loose ends in second review
#
# On the "Materials" sheet:
#  - Remove the Taxon_Local_ID column (was populated with the literal
#    template string ${iNaturalistTaxonId})
#  - Remove the suborder / infraorder / superfamily columns
#  - Fix the scientificNameAuthorship template value from
#    ${summary.Author} to ${summary.authority}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Find the last used column in the header row.
$lastCol = $ws.Cells.Item(1, $ws.Columns.Count).End(-4159).Column

# Collect the column numbers (1-based) of the columns we need to delete,
# identified by their header text in row 1, so the edit is resilient to
# the exact column layout.
$colsToDelete = New-Object System.Collections.ArrayList
for ($c = 1; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item(1, $c).Text
    if ($header -eq "Taxon_Local_ID" -or $header -eq "suborder" -or $header -eq "infraorder" -or $header -eq "superfamily") {
        [void]$colsToDelete.Add($c)
    }
}

# Also find scientificNameAuthorship column so we can fix its value below.
$authorshipCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item(1, $c).Text
    if ($header -eq "scientificNameAuthorship") {
        $authorshipCol = $c
        break
    }
}

# Fix the authorship template value first (column index is unaffected by
# the later column deletions since we do this before deleting anything).
if ($authorshipCol -gt 0) {
    $cell = $ws.Cells.Item(2, $authorshipCol)
    if ($cell.Text -eq "`${summary.Author}") {
        $cell.Value2 = "`${summary.authority}"
    }
}

# Delete the target columns from right to left so earlier column numbers
# stay valid as we go.
$colsToDelete.Sort()
for ($i = $colsToDelete.Count - 1; $i -ge 0; $i--) {
    $colNum = $colsToDelete[$i]
    $ws.Columns.Item($colNum).Delete()
}
